$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new reporting period.
# (Excel automatically repoints the hidden _xlnm._FilterDatabase
# defined name, and the sheet's own autofilter range, to the new name.)
$ws.Name = "Aug 2020 to Sep 2020"
